$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the row-index column style (bold/border/center) from A16 down to the newly appended rows A17:A22
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A22").PasteSpecial(-4122) | Out-Null

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'flowbotallopen'
$ws.Range("C7").Value = 'raw'
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.7272727272727273
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0.5
$ws.Range("K7").Value = 0.9090909090909091
$ws.Range("L7").Value = 0.8
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0.3333333333333333
$ws.Range("O7").Value = 0.8571428571428571
$ws.Range("P7").Value = 0.7516339869281046
$ws.Range("Q7").Value = 0.6666666666666666
$ws.Range("R7").Value = 0.8
$ws.Range("S7").Value = 0.6666666666666666
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0.3333333333333333
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 0.8571428571428571
$ws.Range("X7").Value = 0.75

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'flowbot'
$ws.Range("C8").Value = 'sgp'
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.75
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0.9090909090909091
$ws.Range("L8").Value = 0.8
$ws.Range("M8").Value = 0.3333333333333333
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0.9591836734693876
$ws.Range("P8").Value = 0.915032679738562
$ws.Range("Q8").Value = 0.8888888888888888
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 0.6666666666666666
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 0.3333333333333333
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = 1
$ws.Range("X8").Value = 1

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'dit'
$ws.Range("C9").Value = 'sgp'
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0.5
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.8
$ws.Range("M9").Value = 0.3333333333333333
$ws.Range("N9").Value = 0.9333333333333332
$ws.Range("O9").Value = 0.9183673469387755
$ws.Range("P9").Value = 0.869281045751634
$ws.Range("Q9").Value = 0.7777777777777778
$ws.Range("R9").Value = 0.8
$ws.Range("S9").Value = 0.3333333333333333
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0.6666666666666666
$ws.Range("V9").Value = 1
$ws.Range("W9").Value = 0.7142857142857143
$ws.Range("X9").Value = 1

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'pndit'
$ws.Range("C10").Value = 'sgp'
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0.75
$ws.Range("H10").Value = 0.7272727272727273
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 0.9090909090909091
$ws.Range("L10").Value = 0.8
$ws.Range("M10").Value = 0.6666666666666666
$ws.Range("N10").Value = 0.6
$ws.Range("O10").Value = 0.9183673469387756
$ws.Range("P10").Value = 0.8823529411764706
$ws.Range("Q10").Value = 0.8888888888888888
$ws.Range("R10").Value = 0.8
$ws.Range("S10").Value = 0.3333333333333333
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0.3333333333333333
$ws.Range("V10").Value = 1
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 1

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'hisdit'
$ws.Range("C11").Value = 'sgp'
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.75
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.75
$ws.Range("H11").Value = 0.9090909090909091
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.6
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0.9183673469387756
$ws.Range("P11").Value = 0.954248366013072
$ws.Range("Q11").Value = 0.7777777777777778
$ws.Range("R11").Value = 0.8
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0.6666666666666666
$ws.Range("V11").Value = 1
$ws.Range("W11").Value = 0.8571428571428571
$ws.Range("X11").Value = 1

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 'pnhisdit'
$ws.Range("C12").Value = 'sgp'
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.75
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0.7272727272727273
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0.9090909090909091
$ws.Range("L12").Value = 0.8
$ws.Range("M12").Value = 0.3333333333333333
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0.9795918367346941
$ws.Range("P12").Value = 0.9673202614379085
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 0
$ws.Range("U12").Value = 0.6666666666666666
$ws.Range("V12").Value = 1
$ws.Range("W12").Value = 0.8571428571428571
$ws.Range("X12").Value = 1

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'hisditonly'
$ws.Range("C13").Value = 'sgp'
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.75
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0.25
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.8
$ws.Range("M13").Value = 0.6666666666666666
$ws.Range("N13").Value = 0.7333333333333333
$ws.Range("O13").Value = 0.9183673469387756
$ws.Range("P13").Value = 0.8823529411764706
$ws.Range("Q13").Value = 0.8888888888888888
$ws.Range("R13").Value = 0.8
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0.3333333333333333
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 0.5714285714285714
$ws.Range("X13").Value = 0.25

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 'hisonly'
$ws.Range("C14").Value = 'sgp'
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.25
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0.25
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.5
$ws.Range("K14").Value = 0.9090909090909091
$ws.Range("L14").Value = 0.8
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0.7333333333333333
$ws.Range("O14").Value = 0.9183673469387756
$ws.Range("P14").Value = 0.8431372549019608
$ws.Range("Q14").Value = 0.6666666666666666
$ws.Range("R14").Value = 0.4
$ws.Range("S14").Value = 0.6666666666666666
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0.3333333333333333
$ws.Range("V14").Value = 0.875
$ws.Range("W14").Value = 0.5714285714285714
$ws.Range("X14").Value = 0.5

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 'pndit&pn++'
$ws.Range("C15").Value = 'sgp'
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.5
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.75
$ws.Range("H15").Value = 0.8181818181818182
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.8
$ws.Range("M15").Value = 0.3333333333333333
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0.9387755102040817
$ws.Range("P15").Value = 0.9346405228758168
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 0.3333333333333333
$ws.Range("V15").Value = 1
$ws.Range("W15").Value = 0.8571428571428571
$ws.Range("X15").Value = 1

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 'dit&pn++'
$ws.Range("C16").Value = 'sgp'
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.75
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.75
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.8
$ws.Range("M16").Value = 0.3333333333333333
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0.9183673469387756
$ws.Range("P16").Value = 0.9281045751633988
$ws.Range("Q16").Value = 0.8888888888888888
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 0.3333333333333333
$ws.Range("V16").Value = 1
$ws.Range("W16").Value = 0.8571428571428571
$ws.Range("X16").Value = 1

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 'largedit&pn++'
$ws.Range("C17").Value = 'sgp'
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0.75
$ws.Range("F17").Value = 0.5
$ws.Range("G17").Value = 0.75
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 0.9090909090909091
$ws.Range("L17").Value = 0.8
$ws.Range("M17").Value = 0.3333333333333333
$ws.Range("N17").Value = 0.9333333333333332
$ws.Range("O17").Value = 0.8163265306122449
$ws.Range("P17").Value = 0.9477124183006536
$ws.Range("Q17").Value = 0.7777777777777778
$ws.Range("R17").Value = 0.8
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 0
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = 1
$ws.Range("W17").Value = 0.8571428571428571
$ws.Range("X17").Value = 1

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 'flowbotallopen'
$ws.Range("C18").Value = 'sgp'
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.75
$ws.Range("F18").Value = 0.5
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0.8181818181818182
$ws.Range("L18").Value = 0.8
$ws.Range("M18").Value = 0.3333333333333333
$ws.Range("N18").Value = 0.6
$ws.Range("O18").Value = 0.9387755102040817
$ws.Range("P18").Value = 0.8823529411764706
$ws.Range("Q18").Value = 0.8888888888888888
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0.3333333333333333
$ws.Range("V18").Value = 1
$ws.Range("W18").Value = 0.8571428571428571
$ws.Range("X18").Value = 0.5

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 'pndit&pn++allopen>005'
$ws.Range("C19").Value = 'sgp'
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0.5
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0.9090909090909091
$ws.Range("L19").Value = 0.8
$ws.Range("M19").Value = 0.3333333333333333
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0.9183673469387756
$ws.Range("P19").Value = 0.9738562091503268
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 0.8
$ws.Range("S19").Value = 0.3333333333333333
$ws.Range("T19").Value = 1
$ws.Range("U19").Value = 0.3333333333333333
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 0.8571428571428571
$ws.Range("X19").Value = 1

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 'dit&pn++allopen>01'
$ws.Range("C20").Value = 'sgp'
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.75
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.75
$ws.Range("H20").Value = 0.9090909090909091
$ws.Range("I20").Value = 0.8
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 0.9090909090909091
$ws.Range("L20").Value = 0.6
$ws.Range("M20").Value = 0.3333333333333333
$ws.Range("N20").Value = 0.8
$ws.Range("O20").Value = 0.8571428571428571
$ws.Range("P20").Value = 0.9215686274509803
$ws.Range("Q20").Value = 0.8888888888888888
$ws.Range("R20").Value = 0.8
$ws.Range("S20").Value = 0.6666666666666666
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 0.6666666666666666
$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 0.7142857142857143
$ws.Range("X20").Value = 0.75

# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 'hispnditckpt299'
$ws.Range("C21").Value = 'sgp'
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0.5
$ws.Range("F21").Value = 0.5
$ws.Range("G21").Value = 0.75
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 0.9090909090909091
$ws.Range("L21").Value = 0.8
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0.9591836734693876
$ws.Range("P21").Value = 0.9477124183006538
$ws.Range("Q21").Value = 0.8888888888888888
$ws.Range("R21").Value = 0.6
$ws.Range("S21").Value = 0.6666666666666666
$ws.Range("T21").Value = 1
$ws.Range("U21").Value = 0.3333333333333333
$ws.Range("V21").Value = 1
$ws.Range("W21").Value = 0.8571428571428571
$ws.Range("X21").Value = 1

# Row 22
$ws.Range("A22").Value = 20
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 11
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 11
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 15
$ws.Range("O22").Value = 49
$ws.Range("P22").Value = 153
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 5
$ws.Range("S22").Value = 3
$ws.Range("T22").Value = 1
$ws.Range("U22").Value = 3
$ws.Range("V22").Value = 8
$ws.Range("W22").Value = 7
$ws.Range("X22").Value = 4

